$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the weekly log. It is inserted as
# row 455, which pushes all subsequent rows (old 455-518) down by one
# (new 456-519), matching the diff.
$ws.Rows.Item(455).Insert()

$ws.Cells.Item(455,1).Value = 4
$ws.Cells.Item(455,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(455,3).Value = "Los Lagos"
$ws.Cells.Item(455,4).Value = 45154
$ws.Cells.Item(455,5).Value = 10
$ws.Cells.Item(455,6).Value = 100112045
$ws.Cells.Item(455,7).Value = "Zapallo"
$ws.Cells.Item(455,8).Value = "Paine"
$ws.Cells.Item(455,9).Value = "1a (guarda)"
$ws.Cells.Item(455,10).Value = 250
$ws.Cells.Item(455,11).Value = 600
$ws.Cells.Item(455,12).Value = 600
$ws.Cells.Item(455,13).Value = 600
$ws.Cells.Item(455,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(455,15).Value = "Región de O'Higgins"
$ws.Cells.Item(455,16).Value = 600
$ws.Cells.Item(455,17).Value = 1
$ws.Cells.Item(455,18).Value = "Hortaliza"
